$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (column headers renamed RM1/RM2 -> RC2/RC3) ---
$ws.Range("C1").Value = "RC2"
$ws.Range("D1").Value = "RC3"
$ws.Range("F1").Value = "RC2 posicions"
$ws.Range("G1").Value = "RC3 posicions"

# --- Pre-mark cells whose new text content looks purely numeric, so Excel
#     stores them as text (matching the "posicions" list-of-positions columns),
#     rather than auto-converting to a numeric value. ---
$textForceCells = @("F3", "F4", "F5", "G6", "F8", "G8", "F9", "G9", "G13")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2: ATATATATAT
$ws.Range("A2").Value = "ATATATATAT"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

# Row 3: AAATATATAT
$ws.Range("A3").Value = "AAATATATAT"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "1"
$ws.Range("G3").Value = ""

# Row 4: ATAAATATAT
$ws.Range("A4").Value = "ATAAATATAT"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "3"
$ws.Range("G4").Value = ""

# Row 5: ATATATAAAT
$ws.Range("A5").Value = "ATATATAAAT"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "7"
$ws.Range("G5").Value = ""

# Row 6: TAATATATAT
$ws.Range("A6").Value = "TAATATATAT"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = "1"

# Row 7: ATAAATAAAT
$ws.Range("A7").Value = "ATAAATAAAT"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "3, 7"
$ws.Range("G7").Value = ""

# Row 8: TAATATAAAT
$ws.Range("A8").Value = "TAATATAAAT"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "7"
$ws.Range("G8").Value = "1"

# Row 9: TAAAATATAT
$ws.Range("A9").Value = "TAAAATATAT"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "3"
$ws.Range("G9").Value = "1"

# Row 10: AAAAATATAT
$ws.Range("A10").Value = "AAAAATATAT"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = "1, 3"
$ws.Range("G10").Value = ""

# Row 11: AAATATAAAT
$ws.Range("A11").Value = "AAATATAAAT"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = "1, 7"
$ws.Range("G11").Value = ""

# Row 12: AAAAATAAAT
$ws.Range("A12").Value = "AAAAATAAAT"
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = "1, 3, 7"
$ws.Range("G12").Value = ""

# Row 13: TAAAATAAAT
$ws.Range("A13").Value = "TAAAATAAAT"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = "3, 7"
$ws.Range("G13").Value = "1"

# --- Restore default "Normal" style on the cells we forced to text, so that
#     only the value/type changed and no residual number-format style lingers. ---
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
